$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.676.13'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4916'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2941'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06704'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.901.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07340'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.167'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6677'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.617.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007890'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.130.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.340'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.64%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '191.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.209'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.533'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.932'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.468'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.407'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09155'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.037'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05242'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7429'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.102'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.733'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01822'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.703'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9169'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.069'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +31.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4419'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.913'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9941'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1378'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.552'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.56'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.18%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.056'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05839'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3966'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.61%  '
